# edit.ps1 - Fixed title slides for proj3 lab
#
# Applies:
#  1. Slide 1 title:    "Data 8, Lab 9" -> "Data 8, Project 3 Lab"
#                        (runs: "Data 8" | ", " | "Project 3 Lab", the
#                        latter two colored with the title's accent color)
#  2. Slide 1 subtitle:  "24" + " " + "April 2020" -> single run "24 April 2020"
#  3. Slide 3 title:     "Review: Linear " + "Regression Equation" -> single run
#  4. Slide 6 body:      "Watch " + "year's lecture on privacy " -> single run
#                        (leaves the following hyperlinked "here" run alone)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 1 - Title: "Data 8, Lab " / "9"  ->  "Data 8" / ", " / "Project 3 Lab"
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$titleShape = $slide1.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange

$oldTitle = "Data 8, Lab "
$dataIdx = $titleRange.Text.IndexOf($oldTitle) + 1
# Keep "Data 8" (first 6 chars of "Data 8, Lab ") as-is; replace the
# remainder (", Lab " + the old "9" run) with ", Project 3 Lab".
$tailStart = $dataIdx + 6
$tailLen = $titleRange.Text.Length - ($tailStart - 1)
$tailRange = $titleRange.Characters($tailStart, $tailLen)
$tailRange.Text = ", Project 3 Lab"

# Split that tail into two runs: ", " and "Project 3 Lab"
$commaRange = $titleRange.Characters($tailStart, 2)
$commaRange.Text = ", "
$projRange = $titleRange.Characters($tailStart + 2, 13)
$projRange.Text = "Project 3 Lab"

# Color the two new runs with the title's accent color (C28220)
$accentColor = 0x20 * 65536 + 0x82 * 256 + 0xC2
$commaRange = $titleRange.Characters($tailStart, 2)
$commaRange.Font.Color.RGB = $accentColor
$projRange = $titleRange.Characters($tailStart + 2, 13)
$projRange.Font.Color.RGB = $accentColor

# ---------------------------------------------------------------------
# 2) Slide 1 - Subtitle date line: "24" + " " + "April 2020" -> one run
# ---------------------------------------------------------------------
$subtitleShape = $slide1.Shapes.Item(2)
$subtitleRange = $subtitleShape.TextFrame.TextRange

$dateText = "24 April 2020"
$dateIdx = $subtitleRange.Text.IndexOf("24") + 1
$dateRange = $subtitleRange.Characters($dateIdx, $dateText.Length)
$dateRange.Text = $dateText

# ---------------------------------------------------------------------
# 3) Slide 3 - Title: "Review: Linear " + "Regression Equation" -> one run
# ---------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$s3TitleShape = $slide3.Shapes.Item(1)
$s3TitleRange = $s3TitleShape.TextFrame.TextRange

$reviewText = "Review: Linear Regression Equation"
$reviewRange = $s3TitleRange.Characters(1, $reviewText.Length)
$reviewRange.Text = $reviewText

# ---------------------------------------------------------------------
# 4) Slide 6 - Body: "Watch " + "year's lecture on privacy " -> one run
#    (the hyperlinked "here" run right after must stay untouched)
# ---------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$bodyShape = $slide6.Shapes.Item(2)
$bodyRange = $bodyShape.TextFrame.TextRange

$watchText = "Watch year's lecture on privacy "
$watchIdx = $bodyRange.Text.IndexOf("Watch") + 1
$watchRange = $bodyRange.Characters($watchIdx, $watchText.Length)
$watchRange.Text = $watchText
